$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update C1:C90 from 150 -> 210 and D1:D90 from 2 -> 33
$ws.Range("C1:C90").Value = 210
$ws.Range("D1:D90").Value = 33

# Add new row 91 with A91 = 1
$ws.Range("A91").Value = 1

# Activate the sheet and update the view (topLeftCell + selection)
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 68
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("A1:D90").Select()
